# Update the "tiempo" (D column) measurements in Hoja1 with the new
# timing values from the Radixsort benchmark re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$newValues = @(
    0.0518864,
    0.0370617,
    0.0387722,
    0.0416231,
    0.0456144,
    0.0496057,
    0.0524566,
    0.0553975,
    0.0598689,
    0.0638602,
    0.0678514,
    0.072983,
    0.0758339,
    0.0798252,
    0.0834263,
    0.0860972,
    0.137413,
    0.0963604,
    0.0997815,
    0.103773,
    0.106624,
    0.110045,
    0.113466,
    0.118027,
    0.122589,
    0.12772,
    0.130571,
    0.140264,
    0.18987,
    0.257151,
    0.255441,
    0.25373,
    0.266274,
    0.256011,
    0.261142,
    0.338117,
    0.303906,
    0.282239,
    0.323862,
    0.336976,
    0.293643,
    0.304476,
    0.302195,
    0.302195,
    0.807375,
    0.396275
)

$startRow = 4
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 4).Value = $newValues[$i]
}

# Move the active selection to D50 (below the data range), clearing the
# scrolled viewport back to the top-left of the sheet.
$ws.Range("D50").Select()
